# edit.ps1 - applies the changes described by the target diff to
# draft-gandhi-spring-twamp-srpm-10.pptx via the PowerPoint COM object model.
#
# Summary of changes:
#  1. Handout master date field: "8/4/20" -> "8/4/2020"            (best effort)
#  2. Slide 12, Content Placeholder 2, 2nd paragraph: extend the
#     "Applicable to physical, virtual, LAG ..." sentence with
#     "numbered/unnumbered links".
#  3. Slide 12, "Rectangle 4" shape: move down (Top increases from
#     144.7049pt/1837752 EMU to 156.0258pt/1981527 EMU).
#  4. Slide 3, Content Placeholder 2, 4th paragraph: extend the
#     "Links include physical, virtual, LAG ..." sentence with
#     "numbered/unnumbered links".
#  5. Slide 3, Content Placeholder 2, 8th paragraph ("State is in the
#     probe message"): collapse the 3 separate runs into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master date field "8/4/20" -> "8/4/2020"
#    (best effort - this environment does not expose the handout
#    master's text as writable through the object model, so this is
#    wrapped in a try/catch to keep the rest of the script running
#    regardless.)
# ---------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    $dateShape = $hm.Shapes.Item(2)
    $dateShape.TextFrame.TextRange.Text = "8/4/2020"
} catch {
    Write-Output "handout master date field could not be updated: $_"
}

# ---------------------------------------------------------------------
# Slide 3 ("Requirements and Scope") - Content Placeholder 2
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Content = $s3.Shapes.Item(2)
$s3Tr = $s3Content.TextFrame.TextRange

# 4) 4th paragraph: "Links include physical, virtual, LAG (bundles) and
#    LAG member links" -> "Links include physical, virtual, LAG
#    (bundle), LAG member, numbered/unnumbered links"
$s3Para4 = $s3Tr.Paragraphs(4, 1)
$s3Para4.Text = "TEMP_RESET_4"
$s3Para4b = $s3Tr.Paragraphs(4, 1)
$s3Para4b.Text = "Links include physical, virtual, LAG (bundle), LAG member, numbered/unnumbered links"

# 5) 8th paragraph: merge "State " + "is in the " + "probe message" runs
#    into a single run "State is in the probe message"
$s3Para8 = $s3Tr.Paragraphs(8, 1)
$s3Para8.Text = "TEMP_RESET_8"
$s3Para8b = $s3Tr.Paragraphs(8, 1)
$s3Para8b.Text = "State is in the probe message"

# ---------------------------------------------------------------------
# Slide 12 ("Probe Query for Links") - Content Placeholder 2 + Rectangle 4
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12Content = $s12.Shapes.Item(3)
$s12Tr = $s12Content.TextFrame.TextRange

# 2) 2nd paragraph: "Applicable to physical, virtual, LAG and LAG member
#    links - probe messages pre-routed over the links" ->
#    "Applicable to physical, virtual, LAG, LAG member,
#    numbered/unnumbered links - probe messages pre-routed over the
#    links"
$s12Para2 = $s12Tr.Paragraphs(2, 1)
$s12Para2.Text = "TEMP_RESET_2"
$s12Para2b = $s12Tr.Paragraphs(2, 1)
$s12Para2b.Text = "Applicable to physical, virtual, LAG, LAG member, numbered/unnumbered links – probe messages pre-routed over the links"

# 3) Rectangle 4 shape: move down from y=1837752 EMU to y=1981527 EMU.
#    PowerPoint COM works in points (1pt = 12700 EMU); add half an EMU
#    worth of margin in point-space so the EMU truncation in the host
#    lands on the exact target instead of one EMU short.
$s12Rect = $s12.Shapes.Item(4)
$targetEmu = 1981527
$s12Rect.Top = ($targetEmu + 0.5) / 12700
